$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells we touch in D/E hold plain text (prices/volume %) that must not be
# auto-coerced into numbers by COM, so mark exactly those as Text first.
$textCells = "D2,E2,D3,E3,E4,D5,E5,D6,E6,E7,D8,E8,D9,E9,D10,E10,D11,E11,D12,E12,D13,E13,D14,E14,D15,E15,D16,E16,D17,E17,D18,E18,D19,E19,D20,E20,D21,E21,D22,E22,D23,E23,D24,E24,D25,E25,D26,E26,D27,E27,D28,E28,D29,E29,D30,E30,D31,E31,D32,E32,D33,E33,D34,E34,D35,E35,D36,E36,D37,E37,D38,E38,D39,E39,D40,E40,D41,E41,D42,E42,D43,E43,D44,E44,D45,E45,D46,E46,D47,E47,D48,E48,D49,E49,D50,E50,D51,E51"
foreach ($addr in $textCells.Split(",")) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "33.914.17"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").Value = "1.773.89"
$ws.Range("E3").Value = "  -1.68%  "

$ws.Range("E4").Value = "  +0.51%  "

$ws.Range("D5").Value = "225.59"
$ws.Range("E5").Value = "  -0.68%  "

$ws.Range("D6").Value = "0.546"
$ws.Range("E6").Value = "  +1.64%  "

$ws.Range("E7").Value = "  +0.47%  "

$ws.Range("D8").Value = "30.97"
$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.279"
$ws.Range("E9").Value = "  +0.18%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.0654"
$ws.Range("E10").Value = "  -1.07%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.0928"
$ws.Range("E11").Value = "  +0.12%  "

$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.035.33"
$ws.Range("E12").Value = "  -1.35%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "11.07"
$ws.Range("E13").Value = "  +10.31%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.757.09"
$ws.Range("E14").Value = "  -2.63%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.624"
$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "33.952.09"
$ws.Range("E16").Value = "  +0.57%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "4.19"
$ws.Range("E17").Value = "  -1.03%  "

$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "68.74"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "250.56"
$ws.Range("E19").Value = "  -1.66%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0735"
$ws.Range("E20").Value = "  -0.60%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "10.31"
$ws.Range("E22").Value = "  -0.79%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "4.18"
$ws.Range("E23").Value = "  -2.65%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "2.14"
$ws.Range("E24").Value = "  -1.94%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "156.14"
$ws.Range("E25").Value = "  -0.25%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "16.31"
$ws.Range("E26").Value = "  -0.38%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "6.94"
$ws.Range("E27").Value = "  -0.96%  "

$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "0.113"
$ws.Range("E28").Value = "  -1.15%  "

$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.51%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "3.74"
$ws.Range("E30").Value = "  -1.74%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.0512"
$ws.Range("E31").Value = "  +1.16%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").Value = "  -0.12%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "3.55"
$ws.Range("E33").Value = "  +2.22%  "

$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "1.81"
$ws.Range("E34").Value = "  +1.85%  "

$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "1.443.22"
$ws.Range("E35").Value = "  -5.76%  "

$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "1.06"
$ws.Range("E36").Value = "  -0.99%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.623"
$ws.Range("E37").Value = "  +1.99%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0186"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.84"
$ws.Range("E39").Value = "  +1.77%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "82.12"
$ws.Range("E40").Value = "  -1.45%  "

$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "2.35"
$ws.Range("E41").Value = "  +1.38%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "0.883"
$ws.Range("E42").Value = "  -2.02%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "2.05"
$ws.Range("E43").Value = "  -2.66%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "0.0508"
$ws.Range("E44").Value = "  -2.38%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "1.06"
$ws.Range("E45").Value = "  -0.93%  "

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.932.98"
$ws.Range("E46").Value = "  -1.03%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "5.75"
$ws.Range("E47").Value = "  +2.06%  "

$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.53%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "11.71"
$ws.Range("E49").Value = "  +5.37%  "

$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "49.70"
$ws.Range("E50").Value = "  -4.28%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "97.07"
$ws.Range("E51").Value = "  +2.18%  "
